$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''246.59'
$ws.Range("G2").Value = '''20'
$ws.Range("D3").Value = '''29.60'
$ws.Range("E3").Value = '''9.14%'
$ws.Range("G3").Value = '''20'
$ws.Range("D4").Value = '''5.169'
$ws.Range("E4").Value = '''1.19%'
$ws.Range("G4").Value = '''20'
$ws.Range("E5").Value = '''0.46%'
$ws.Range("G5").Value = '''20'
$ws.Range("D6").Value = '''6.584'
$ws.Range("E6").Value = '''1.19%'
$ws.Range("G6").Value = '''20'
$ws.Range("D7").Value = '''0.8572'
$ws.Range("E7").Value = '''4.61%'
$ws.Range("G7").Value = '''20'
$ws.Range("D8").Value = '''0.8696'
$ws.Range("E8").Value = '''0.69%'
$ws.Range("G8").Value = '''20'
$ws.Range("E9").Value = '''2.69%'
$ws.Range("G9").Value = '''20'
$ws.Range("D10").Value = '''0.07076'
$ws.Range("E10").Value = '''1.95%'
$ws.Range("G10").Value = '''20'
$ws.Range("D11").Value = '''0.02925'
$ws.Range("E11").Value = '''2.91%'
$ws.Range("G11").Value = '''20'
$ws.Range("D12").Value = '''0.09382'
$ws.Range("E12").Value = '''-0.13%'
$ws.Range("G12").Value = '''20'
$ws.Range("D13").Value = '''0.001519'
$ws.Range("E13").Value = '''0.09%'
$ws.Range("G13").Value = '''20'
$ws.Range("D14").Value = '''0.04139'
$ws.Range("E14").Value = '''2.74%'
$ws.Range("G14").Value = '''20'
$ws.Range("D15").Value = '''0.0006021'
$ws.Range("E15").Value = '''-94.04%'
$ws.Range("G15").Value = '''20'
$ws.Range("D16").Value = '''0.006173'
$ws.Range("E16").Value = '''0.77%'
$ws.Range("G16").Value = '''20'
$ws.Range("E17").Value = '''3,765.12%'
$ws.Range("G17").Value = '''20'
$ws.Range("D18").Value = '''3.491'
$ws.Range("E18").Value = '''-0.58%'
$ws.Range("G18").Value = '''20'
$ws.Range("D19").Value = '''3.100'
$ws.Range("E19").Value = '''3.01%'
$ws.Range("G19").Value = '''20'
$ws.Range("D20").Value = '''2.278'
$ws.Range("E20").Value = '''-1.75%'
$ws.Range("G20").Value = '''20'
$ws.Range("D21").Value = '''0.3181'
$ws.Range("E21").Value = '''0.52%'
$ws.Range("G21").Value = '''20'
$ws.Range("D22").Value = '''0.03391'
$ws.Range("E22").Value = '''5.82%'
$ws.Range("G22").Value = '''20'
$ws.Range("E23").Value = '''0.04%'
$ws.Range("G23").Value = '''20'
$ws.Range("D24").Value = '''3.468'
$ws.Range("E24").Value = '''-3.23%'
$ws.Range("G24").Value = '''20'
$ws.Range("D25").Value = '''0.1380'
$ws.Range("E25").Value = '''0.47%'
$ws.Range("G25").Value = '''20'
$ws.Range("D26").Value = '''0.005006'
$ws.Range("E26").Value = '''12.02%'
$ws.Range("G26").Value = '''20'
$ws.Range("D27").Value = '''0.001223'
$ws.Range("E27").Value = '''0.34%'
$ws.Range("G27").Value = '''20'
$ws.Range("G28").Value = '''20'
$ws.Range("G29").Value = '''20'
$ws.Range("G30").Value = '''20'
$ws.Range("G31").Value = '''20'
$ws.Range("G32").Value = '''20'
$ws.Range("G33").Value = '''20'
$ws.Range("G34").Value = '''20'
$ws.Range("G35").Value = '''20'
$ws.Range("G36").Value = '''20'
$ws.Range("G37").Value = '''20'
$ws.Range("G38").Value = '''20'
$ws.Range("G39").Value = '''20'
$ws.Range("D40").Value = '''0.03752'
$ws.Range("E40").Value = '''0.64%'
$ws.Range("G40").Value = '''20'
$ws.Range("D41").Value = '''0.005758'
$ws.Range("E41").Value = '''-3.38%'
$ws.Range("G41").Value = '''20'
$ws.Range("D42").Value = '''0.1070'
$ws.Range("E42").Value = '''1.13%'
$ws.Range("G42").Value = '''20'
$ws.Range("D43").Value = '''0.002383'
$ws.Range("E43").Value = '''3.61%'
$ws.Range("G43").Value = '''20'
$ws.Range("D44").Value = '''0.008481'
$ws.Range("E44").Value = '''-11.12%'
$ws.Range("G44").Value = '''20'
$ws.Range("D45").Value = '''0.00005260'
$ws.Range("E45").Value = '''2.19%'
$ws.Range("G45").Value = '''20'
$ws.Range("E46").Value = '''-0.02%'
$ws.Range("G46").Value = '''20'
$ws.Range("D47").Value = '''0.06471'
$ws.Range("E47").Value = '''-35.92%'
$ws.Range("G47").Value = '''20'
$ws.Range("D48").Value = '''0.002532'
$ws.Range("E48").Value = '''1.26%'
$ws.Range("G48").Value = '''20'
$ws.Range("E49").Value = '''-0.02%'
$ws.Range("G49").Value = '''20'
$ws.Range("E50").Value = '''-0.02%'
$ws.Range("G50").Value = '''20'
$ws.Range("G51").Value = '''20'
